$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force text-typed values (avoids Excel auto-numeric
# coercion on number-looking strings) without altering the target cell style.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

function Set-TextValue($cell, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163) | Out-Null
}

# Row 2
Set-TextValue $ws.Range("D2") "43.355.27"
$ws.Range("E2").Value = "  +2.70%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.305.21"
$ws.Range("E3").Value = "  +1.67%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
Set-TextValue $ws.Range("D5") "311.19"
$ws.Range("E5").Value = "  +1.51%  "

# Row 6
Set-TextValue $ws.Range("D6") "102.47"
$ws.Range("E6").Value = "  +5.47%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.531"
$ws.Range("E7").Value = "  +1.07%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.529"

# Row 10
Set-TextValue $ws.Range("D10") "35.66"
$ws.Range("E10").Value = "  +1.49%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0811"
$ws.Range("E11").Value = "  +2.68%  "

# Row 12
$ws.Range("E12").Value = "  -1.01%  "

# Row 13
Set-TextValue $ws.Range("D13") "6.97"
$ws.Range("E13").Value = "  +1.43%  "

# Row 14
Set-TextValue $ws.Range("D14") "2.663.67"
$ws.Range("E14").Value = "  +1.65%  "

# Row 15
Set-TextValue $ws.Range("D15") "14.99"
$ws.Range("E15").Value = "  +2.05%  "

# Row 16
Set-TextValue $ws.Range("D16") "2.281.77"
$ws.Range("E16").Value = "  +1.06%  "

# Row 17
$ws.Range("E17").Value = "  +1.91%  "

# Row 18
Set-TextValue $ws.Range("D18") "43.264.71"
$ws.Range("E18").Value = "  +2.71%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.29"
$ws.Range("E19").Value = "  +0.26%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.0₃0933"
$ws.Range("E20").Value = "  +3.08%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.17"

# Row 22
Set-TextValue $ws.Range("D22") "68.06"
$ws.Range("E22").Value = "  +0.38%  "

# Row 23
Set-TextValue $ws.Range("D23") "241.28"
$ws.Range("E23").Value = "  +1.75%  "

# Row 24
$ws.Range("E24").Value = "  +1.29%  "

# Row 25
$ws.Range("E25").Value = "  +1.78%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("E27").Value = "  -1.65%  "

# Row 28
Set-TextValue $ws.Range("D28") "24.74"
$ws.Range("E28").Value = "  +5.13%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D29") "36.70"
$ws.Range("E29").Value = "  -2.72%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "2.20"
$ws.Range("E30").Value = "  +3.81%  "

# Row 31
Set-TextValue $ws.Range("D31") "9.63"
$ws.Range("E31").Value = "  +0.50%  "

# Row 32
Set-TextValue $ws.Range("D32") "168.76"
$ws.Range("E32").Value = "  +3.84%  "

# Row 33
Set-TextValue $ws.Range("D33") "5.27"
$ws.Range("E33").Value = "  +0.46%  "

# Row 34
$ws.Range("E34").Value = "  +0.01%  "

# Row 35
Set-TextValue $ws.Range("D35") "2.53"
$ws.Range("E35").Value = "  +6.63%  "

# Row 36
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D36") "17.72"
$ws.Range("E36").Value = "  +0.65%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D37") "0.0742"
$ws.Range("E37").Value = "  +0.77%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.06"
$ws.Range("E38").Value = "  -2.72%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.89"
$ws.Range("E39").Value = "  +3.83%  "

# Row 40
$ws.Range("E40").Value = "  +1.84%  "

# Row 41
$ws.Range("E41").Value = "  +1.14%  "

# Row 42
Set-TextValue $ws.Range("D42") "4.35"
$ws.Range("E42").Value = "  +6.03%  "

# Row 43
$ws.Range("E43").Value = "  -1.29%  "

# Row 44
Set-TextValue $ws.Range("D44") "19.42"
$ws.Range("E44").Value = "  +2.20%  "

# Row 45
$ws.Range("E45").Value = "  +2.76%  "

# Row 46
Set-TextValue $ws.Range("D46") "1.965.64"
$ws.Range("E46").Value = "  +0.77%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.99"
$ws.Range("E47").Value = "  +2.46%  "

# Row 48
Set-TextValue $ws.Range("D48") "9.88"
$ws.Range("E48").Value = "  -1.30%  "

# Row 49
Set-TextValue $ws.Range("D49") "55.39"
$ws.Range("E49").Value = "  +2.71%  "

# Row 50
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D50") "2.90"
$ws.Range("E50").Value = "  +1.02%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D51") "1.58"
$ws.Range("E51").Value = "  +7.25%  "

$scratch.Clear()
$ws.Range("A1").Select() | Out-Null
